$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "want to go" counts in column F
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 1970
$wsExhibition.Range("F4").Value = 843
$wsExhibition.Range("F5").Value = 1007

# Sheet "全部类型" (All types) - update "want to go" counts in column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1970
$wsAll.Range("F5").Value = 843
$wsAll.Range("F6").Value = 1007
